$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 425
$ws.Range("I8").Value = 425
$ws.Range("K8").Value = 1275
$ws.Range("M8").Value = -1136
$ws.Range("H18").Value = 763.5714
$ws.Range("I18").Value = 724.1667
$ws.Range("K18").Value = 724.1667
$ws.Range("M18").Value = -440.1667
$ws.Range("H19").Value = 1588.1
$ws.Range("I19").Value = 1849.6
$ws.Range("J19").Value = 1326.6
$ws.Range("K19").Value = 1849.6
$ws.Range("L19").Value = 1326.6
$ws.Range("M19").Value = -1674.6
$ws.Range("N19").Value = -1676.6
$ws.Range("H32").Value = 2445.3572
$ws.Range("J32").Value = 2745.6
$ws.Range("L32").Value = 2745.6
$ws.Range("N32").Value = -3397.6
$ws.Range("H51").Value = 8749.700000000001
$ws.Range("I51").Value = 11261
$ws.Range("J51").Value = 6238.4
$ws.Range("K51").Value = 11261
$ws.Range("L51").Value = 6238.4
$ws.Range("M51").Value = -10777
$ws.Range("N51").Value = -7206.4
$ws.Range("H64").Value = 4988
$ws.Range("I64").Value = 4983.5
$ws.Range("J64").Value = 4992.5
$ws.Range("K64").Value = 4983.5
$ws.Range("L64").Value = 4992.5
$ws.Range("M64").Value = -4735.5
$ws.Range("N64").Value = -5488.5
$ws.Range("H67").Value = 4988
$ws.Range("I67").Value = 4983.5
$ws.Range("J67").Value = 4992.5
$ws.Range("K67").Value = 4983.5
$ws.Range("L67").Value = 4992.5
$ws.Range("M67").Value = -4125.5
$ws.Range("N67").Value = -6708.5
$ws.Range("H69").Value = 15124.25
$ws.Range("I69").Value = 10498.5
$ws.Range("K69").Value = 31495.5
$ws.Range("M69").Value = -30621.5
$ws.Range("H72").Value = 15124.25
$ws.Range("I72").Value = 10498.5
$ws.Range("K72").Value = 94486.5
$ws.Range("M72").Value = -90118.5
$ws.Range("H76").Value = 4999.5
$ws.Range("I76").Value = 4999.5
$ws.Range("K76").Value = 4999.5
$ws.Range("M76").Value = -4684.5
$ws.Range("H79").Value = 4999.5
$ws.Range("I79").Value = 4999.5
$ws.Range("K79").Value = 4999.5
$ws.Range("M79").Value = -3907.5
$ws.Range("H80").Value = 35998.06
$ws.Range("I80").Value = 86053.86
$ws.Range("J80").Value = 959
$ws.Range("K80").Value = 258161.58
$ws.Range("L80").Value = 2877
$ws.Range("M80").Value = -257163.58
$ws.Range("N80").Value = -4873
$ws.Range("H83").Value = 35998.06
$ws.Range("I83").Value = 86053.86
$ws.Range("J83").Value = 959
$ws.Range("K83").Value = 774484.74
$ws.Range("L83").Value = 8631
$ws.Range("M83").Value = -769492.74
$ws.Range("N83").Value = -18615
$ws.Range("H116").Value = 5129.857
$ws.Range("I116").Value = 4896.04
$ws.Range("K116").Value = 4896.04
$ws.Range("M116").Value = -1454.04
$ws.Range("H125").Value = 1418.8
$ws.Range("I125").Value = 562.5
$ws.Range("J125").Value = 1989.6666
$ws.Range("K125").Value = 5062.5
$ws.Range("L125").Value = 17906.9994
$ws.Range("M125").Value = -2602.5
$ws.Range("N125").Value = -22826.9994
$ws.Range("H137").Value = 1760.2222
$ws.Range("I137").Value = 950.1667
$ws.Range("J137").Value = 3380.3333
$ws.Range("K137").Value = 2850.5001
$ws.Range("L137").Value = 10140.9999
$ws.Range("M137").Value = -300.5001000000002
$ws.Range("N137").Value = -15240.9999
$ws.Range("H138").Value = 4154.381
$ws.Range("I138").Value = 3485.2942
$ws.Range("J138").Value = 6998
$ws.Range("K138").Value = 10455.8826
$ws.Range("L138").Value = 20994
$ws.Range("M138").Value = -5315.882599999999
$ws.Range("N138").Value = -31274
$ws.Range("H141").Value = 6721.25
$ws.Range("I141").Value = 6721.25
$ws.Range("K141").Value = 20163.75
$ws.Range("M141").Value = -14983.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31658.816
$ws.Range("I32").Value = 35737.387
$ws.Range("J32").Value = 13596.571
$ws.Range("K32").Value = 35737.387
$ws.Range("L32").Value = 13596.571
$ws.Range("M32").Value = -35450.387
$ws.Range("N32").Value = -14170.571
$ws.Range("H41").Value = 2009.6666
$ws.Range("I41").Value = 1264.5
$ws.Range("K41").Value = 1264.5
$ws.Range("M41").Value = -850.5
$ws.Range("H45").Value = 2283.1875
$ws.Range("I45").Value = 1034
$ws.Range("J45").Value = 3889.2856
$ws.Range("K45").Value = 1034
$ws.Range("L45").Value = 3889.2856
$ws.Range("M45").Value = -657
$ws.Range("N45").Value = -4643.2856
$ws.Range("H61").Value = 2986.875
$ws.Range("I61").Value = 2986.875
$ws.Range("K61").Value = 2986.875
$ws.Range("M61").Value = -2774.875
$ws.Range("H74").Value = 35399.723
$ws.Range("I74").Value = 36592.57
$ws.Range("K74").Value = 36592.57
$ws.Range("M74").Value = -35718.57
$ws.Range("H77").Value = 35399.723
$ws.Range("I77").Value = 36592.57
$ws.Range("K77").Value = 182962.85
$ws.Range("M77").Value = -178594.85
$ws.Range("H97").Value = 5311.393
$ws.Range("I97").Value = 5582.773
$ws.Range("J97").Value = 4316.3335
$ws.Range("K97").Value = 5582.773
$ws.Range("L97").Value = 4316.3335
$ws.Range("M97").Value = -5086.773
$ws.Range("N97").Value = -5308.3335
$ws.Range("H102").Value = 2564
$ws.Range("I102").Value = 2355.0625
$ws.Range("K102").Value = 2355.0625
$ws.Range("M102").Value = -733.0625
$ws.Range("H122").Value = 1377.3572
$ws.Range("I122").Value = 1450.25
$ws.Range("J122").Value = 940
$ws.Range("K122").Value = 4350.75
$ws.Range("L122").Value = 2820
$ws.Range("M122").Value = -1900.75
$ws.Range("N122").Value = -7720
$ws.Range("H132").Value = 93302
$ws.Range("I132").Value = 93302
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 279906
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -277376
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2986.875
$ws.Range("I136").Value = 2986.875
$ws.Range("K136").Value = 8960.625
$ws.Range("M136").Value = -6410.625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1004
$ws.Range("I8").Value = 1004
$ws.Range("K8").Value = 1004
$ws.Range("M8").Value = -864
$ws.Range("H20").Value = 3410.35
$ws.Range("J20").Value = 3099.6667
$ws.Range("L20").Value = 3099.6667
$ws.Range("N20").Value = -3593.6667
$ws.Range("H27").Value = 99997.5
$ws.Range("J27").Value = 99997.5
$ws.Range("L27").Value = 99997.5
$ws.Range("N27").Value = -100381.5
$ws.Range("H60").Value = 69932.336
$ws.Range("I60").Value = 69898
$ws.Range("K60").Value = 69898
$ws.Range("M60").Value = -69299
$ws.Range("H99").Value = 56577.316
$ws.Range("I99").Value = 86846.586
$ws.Range("K99").Value = 86846.586
$ws.Range("M99").Value = -85348.586
$ws.Range("H105").Value = 2992
$ws.Range("I105").Value = 3137.6
$ws.Range("J105").Value = 2409.6
$ws.Range("K105").Value = 3137.6
$ws.Range("L105").Value = 2409.6
$ws.Range("M105").Value = -1390.6
$ws.Range("N105").Value = -5903.6
$ws.Range("H107").Value = 1090.5
$ws.Range("I107").Value = 1090.5
$ws.Range("K107").Value = 1090.5
$ws.Range("M107").Value = 829.5
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("N122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("N123").Value = 0
$ws.Range("N123").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 721.6
$ws.Range("I2").Value = 852
$ws.Range("K2").Value = 852
$ws.Range("M2").Value = -739
$ws.Range("H16").Value = 582.55554
$ws.Range("I16").Value = 599.36
$ws.Range("J16").Value = 372.5
$ws.Range("K16").Value = 599.36
$ws.Range("L16").Value = 372.5
$ws.Range("M16").Value = -312.36
$ws.Range("N16").Value = -946.5
$ws.Range("H31").Value = 2554.3684
$ws.Range("I31").Value = 1823
$ws.Range("J31").Value = 3560
$ws.Range("K31").Value = 1823
$ws.Range("L31").Value = 3560
$ws.Range("M31").Value = -1528
$ws.Range("N31").Value = -4150
$ws.Range("H33").Value = 16500
$ws.Range("I33").Value = 1000
$ws.Range("K33").Value = 1000
$ws.Range("M33").Value = -621
$ws.Range("H34").Value = 2554.3684
$ws.Range("I34").Value = 1823
$ws.Range("J34").Value = 3560
$ws.Range("K34").Value = 1823
$ws.Range("L34").Value = 3560
$ws.Range("M34").Value = -1621
$ws.Range("N34").Value = -3964
$ws.Range("H86").Value = 61336.57
$ws.Range("H89").Value = 61336.57
$ws.Range("H105").Value = 1259.8462
$ws.Range("I105").Value = 1259.8462
$ws.Range("K105").Value = 1259.8462
$ws.Range("M105").Value = 487.1538
$ws.Range("H107").Value = 445.1579
$ws.Range("I107").Value = 424
$ws.Range("J107").Value = 625
$ws.Range("K107").Value = 424
$ws.Range("L107").Value = 625
$ws.Range("M107").Value = 1496
$ws.Range("N107").Value = -4465
$ws.Range("H113").Value = 582.55554
$ws.Range("I113").Value = 599.36
$ws.Range("J113").Value = 372.5
$ws.Range("K113").Value = 599.36
$ws.Range("L113").Value = 372.5
$ws.Range("M113").Value = 1570.64
$ws.Range("N113").Value = -4712.5
$ws.Range("H122").Value = 1191.421
$ws.Range("J122").Value = 1243.3334
$ws.Range("L122").Value = 3730.0002
$ws.Range("N122").Value = -8630.0002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 64625
$ws.Range("J37").Value = 64625
$ws.Range("L37").Value = 193875
$ws.Range("N37").Value = -194099
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("N75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H93").Value = 10113.272
$ws.Range("J93").Value = 12222.223
$ws.Range("L93").Value = 36666.669
$ws.Range("N93").Value = -40410.669
$ws.Range("H96").Value = 10025
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 1058
$ws.Range("J122").Value = 1058
$ws.Range("L122").Value = 9522
$ws.Range("N122").Value = -14422
$ws.Range("H131").Value = 2229632.8
$ws.Range("I131").Value = 12864.556
$ws.Range("J131").Value = 2783824.8
$ws.Range("K131").Value = 38593.66800000001
$ws.Range("L131").Value = 8351474.399999999
$ws.Range("M131").Value = -33553.66800000001
$ws.Range("N131").Value = -8361554.399999999
$ws.Range("H132").Value = 1166
$ws.Range("I132").Value = 1166
$ws.Range("K132").Value = 10494
$ws.Range("M132").Value = -7964
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4701.4
$ws.Range("I43").Value = 4701.4
$ws.Range("K43").Value = 4701.4
$ws.Range("M43").Value = -4550.4
$ws.Range("H46").Value = 17220.5
$ws.Range("I46").Value = 17220.5
$ws.Range("K46").Value = 17220.5
$ws.Range("M46").Value = -17064.5
$ws.Range("H57").Value = 16178.667
$ws.Range("I57").Value = 10268.333
$ws.Range("J57").Value = 27999.334
$ws.Range("K57").Value = 10268.333
$ws.Range("L57").Value = 27999.334
$ws.Range("M57").Value = -9448.333000000001
$ws.Range("N57").Value = -29639.334
$ws.Range("H80").Value = 3040.8333
$ws.Range("I80").Value = 2563
$ws.Range("J80").Value = 3518.6667
$ws.Range("K80").Value = 2563
$ws.Range("L80").Value = 3518.6667
$ws.Range("M80").Value = -1565
$ws.Range("N80").Value = -5514.6667
$ws.Range("H83").Value = 3040.8333
$ws.Range("I83").Value = 2563
$ws.Range("J83").Value = 3518.6667
$ws.Range("K83").Value = 12815
$ws.Range("L83").Value = 17593.3335
$ws.Range("M83").Value = -7823
$ws.Range("N83").Value = -27577.3335
$ws.Range("H102").Value = 2522.6924
$ws.Range("I102").Value = 1654.6666
$ws.Range("J102").Value = 4475.75
$ws.Range("K102").Value = 1654.6666
$ws.Range("L102").Value = 4475.75
$ws.Range("M102").Value = -32.66660000000002
$ws.Range("N102").Value = -7719.75
$ws.Range("H122").Value = 1867.6666
$ws.Range("I122").Value = 1907
$ws.Range("J122").Value = 1199
$ws.Range("K122").Value = 5721
$ws.Range("L122").Value = 3597
$ws.Range("M122").Value = -3271
$ws.Range("N122").Value = -8497
$ws.Range("H126").Value = 7471.643
$ws.Range("I126").Value = 6319.6665
$ws.Range("K126").Value = 18958.9995
$ws.Range("M126").Value = -16488.9995
$ws.Range("H132").Value = 34683.355
$ws.Range("I132").Value = 42198.12
$ws.Range("J132").Value = 3371.8333
$ws.Range("K132").Value = 126594.36
$ws.Range("L132").Value = 10115.4999
$ws.Range("M132").Value = -124064.36
$ws.Range("N132").Value = -15175.4999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4175.9443
$ws.Range("I40").Value = 3462
$ws.Range("J40").Value = 6674.75
$ws.Range("K40").Value = 3462
$ws.Range("L40").Value = 6674.75
$ws.Range("M40").Value = -3326
$ws.Range("N40").Value = -6946.75
$ws.Range("H61").Value = 4091.8823
$ws.Range("I61").Value = 3972.625
$ws.Range("K61").Value = 3972.625
$ws.Range("M61").Value = -3770.625
$ws.Range("H113").Value = 4091.8823
$ws.Range("I113").Value = 3972.625
$ws.Range("K113").Value = 3972.625
$ws.Range("M113").Value = -1802.625
$ws.Range("H122").Value = 3896.8333
$ws.Range("I122").Value = 2799.6667
$ws.Range("K122").Value = 8399.000100000001
$ws.Range("M122").Value = -5949.000100000001
$ws.Range("H124").Value = 150214.5
$ws.Range("J124").Value = 150214.5
$ws.Range("L124").Value = 150214.5
$ws.Range("N124").Value = -160034.5
$ws.Range("H132").Value = 84767.13
$ws.Range("I132").Value = 112227.91
$ws.Range("K132").Value = 336683.73
$ws.Range("M132").Value = -334153.73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 16999.5
$ws.Range("I10").Value = 15999
$ws.Range("J10").Value = 18000
$ws.Range("K10").Value = 15999
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = -15830
$ws.Range("N10").Value = -18338
$ws.Range("H32").Value = 10433
$ws.Range("I32").Value = 10433
$ws.Range("K32").Value = 10433
$ws.Range("M32").Value = -10116
$ws.Range("H49").Value = 24050
$ws.Range("J49").Value = 24050
$ws.Range("L49").Value = 24050
$ws.Range("N49").Value = -24510
$ws.Range("H122").Value = 1689.091
$ws.Range("I122").Value = 1734.4736
$ws.Range("K122").Value = 5203.4208
$ws.Range("M122").Value = -2753.4208
$ws.Range("H132").Value = 31820.846
$ws.Range("I132").Value = 33217.34
$ws.Range("K132").Value = 99652.01999999999
$ws.Range("M132").Value = -97122.01999999999
$ws.Range("H136").Value = 3056.7917
$ws.Range("I136").Value = 2728.2856
$ws.Range("K136").Value = 8184.8568
$ws.Range("M136").Value = -5634.8568
